$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: new log entry - date (formatted like the existing entries), task, associated file
$ws.Range("A2").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = 44497

$ws.Range("B16").Value = "Added Shrew #4 with labels - very careful in labeling, excluded blurry images"

# Row 17: continuation task note (no date / associated file, matching row 3's pattern)
$ws.Range("B17").Value = "Trained at 2000 iter pcutoff = 0.9, MSE ~5"

$ws.Range("C16").Value = "labeled_data"

$ws.Range("C16").Select()
